$wb = $excel.ActiveWorkbook

# Locate source sheet ("moodle1") that will be duplicated as the new sheet.
$src = $wb.Worksheets.Item("moodle1")

# Add a brand-new worksheet positioned right after "moodle1" and name it "moodle2".
$new = $wb.Worksheets.Add($null, $src)
$new.Name = "moodle2"

# Copy every used cell (values) from "moodle1" into "moodle2".
$usedRange = $src.UsedRange
$rows = $usedRange.Rows.Count
$cols = $usedRange.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $val = $src.Cells.Item($r, $c).Value2
        if ($null -ne $val) {
            $new.Cells.Item($r, $c).Value2 = $val
        }
    }
}

$new.Range("A1").Select()
